$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column (Price) values: force Text type first so Excel's COM layer
# doesn't silently reinterpret numeric-looking strings as floating point
# numbers (which would corrupt exact text like '560.19' into
# '560.19000000000005', or strip formatting like '64.029.50').
$dCells = @('D2', 'D3', 'D5', 'D6', 'D8', 'D9', 'D10', 'D11', 'D12', 'D14', 'D15', 'D16', 'D17', 'D19', 'D20', 'D21', 'D22', 'D23', 'D25', 'D27', 'D28', 'D31', 'D34', 'D35', 'D36', 'D38', 'D40', 'D41', 'D45', 'D49', 'D50')
foreach ($c in $dCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '63.977.39'
$ws.Range('D3').Value = '3.061.72'
$ws.Range('D5').Value = '560.19'
$ws.Range('D6').Value = '143.79'
$ws.Range('D8').Value = '3.060.15'
$ws.Range('D9').Value = '0.514'
$ws.Range('D10').Value = '0.154'
$ws.Range('D11').Value = '6.19'
$ws.Range('D12').Value = '0.497'
$ws.Range('D14').Value = '35.92'
$ws.Range('D15').Value = '3.563.02'
$ws.Range('D16').Value = '64.035.91'
$ws.Range('D17').Value = '3.060.09'
$ws.Range('D19').Value = '6.82'
$ws.Range('D20').Value = '477.85'
$ws.Range('D21').Value = '14.06'
$ws.Range('D22').Value = '0.686'
$ws.Range('D23').Value = '14.48'
$ws.Range('D25').Value = '82.53'
$ws.Range('D27').Value = '2.81'
$ws.Range('D28').Value = '8.15'
$ws.Range('D31').Value = '26.34'
$ws.Range('D34').Value = '5.78'
$ws.Range('D35').Value = '6.25'
$ws.Range('D36').Value = '54.68'
$ws.Range('D38').Value = '453.10'
$ws.Range('D40').Value = '2.83'
$ws.Range('D41').Value = '3.023.00'
$ws.Range('D45').Value = '27.95'
$ws.Range('D49').Value = '119.11'
$ws.Range('D50').Value = '0.0₃0518'

# Restore the cells' style to the sheet default (Normal) now that the
# text value is committed, so no stray number-format style lingers.
foreach ($c in $dCells) {
    $ws.Range($c).Style = "Normal"
}

# E-column (Volume 1h) values already contain non-numeric characters
# ('%', leading sign, padding spaces) so Excel keeps them as text as-is.
$ws.Range('E2').Value = '  -0.82%  '
$ws.Range('E3').Value = '  -0.58%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E5').Value = '  +0.98%  '
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -0.49%  '
$ws.Range('E9').Value = '  +3.84%  '
$ws.Range('E10').Value = '  +1.67%  '
$ws.Range('E11').Value = '  -10.38%  '
$ws.Range('E12').Value = '  +9.11%  '
$ws.Range('E13').Value = '  +2.45%  '
$ws.Range('E14').Value = '  +1.94%  '
$ws.Range('E15').Value = '  -0.42%  '
$ws.Range('E16').Value = '  -0.69%  '
$ws.Range('E17').Value = '  -0.75%  '
$ws.Range('E18').Value = '  +1.01%  '
$ws.Range('E19').Value = '  +1.53%  '
$ws.Range('E20').Value = '  -0.59%  '
$ws.Range('E22').Value = '  +2.35%  '
$ws.Range('E23').Value = '  +10.29%  '
$ws.Range('E24').Value = '  +0.82%  '
$ws.Range('E25').Value = '  +2.17%  '
$ws.Range('E26').Value = '  -0.71%  '
$ws.Range('E27').Value = '  +0.50%  '
$ws.Range('E28').Value = '  +3.22%  '
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('E31').Value = '  +1.08%  '
$ws.Range('E32').Value = '  -0.48%  '
$ws.Range('E33').Value = '  +1.82%  '
$ws.Range('E34').Value = '  +1.31%  '
$ws.Range('E35').Value = '  +3.04%  '
$ws.Range('E36').Value = '  -1.09%  '
$ws.Range('E37').Value = '  +1.45%  '
$ws.Range('E38').Value = '  -2.42%  '
$ws.Range('E39').Value = '  -1.01%  '
$ws.Range('E40').Value = '  +5.05%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('E42').Value = '  +0.60%  '
$ws.Range('E43').Value = '  -1.46%  '
$ws.Range('E44').Value = '  +4.36%  '
$ws.Range('E45').Value = '  +0.49%  '
$ws.Range('E46').Value = '  +9.13%  '
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('E48').Value = '  +1.82%  '
$ws.Range('E49').Value = '  +1.39%  '
$ws.Range('E50').Value = '  +0.85%  '
$ws.Range('E51').Value = '  +2.73%  '
